$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-15 09:00:37"
$wsZh.Range("D3").Value = "2016-02-15 09:00:37"
$wsZh.Range("G2").Value = "2016-02-15 09:01:55"
$wsZh.Range("G3").Value = "2016-02-15 09:01:55"

# "de-de" sheet: Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-15 09:00:56"
$wsDe.Range("D3").Value = "2016-02-15 09:00:56"
$wsDe.Range("G2").Value = "2016-02-15 09:02:23"
$wsDe.Range("G3").Value = "2016-02-15 09:02:23"
